# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback (target files + handback xlf files) has been generated/received
# for both the zh-cn and de-de locales, and that the status moved from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdFile     = "4f1cbecd-55f8-4daa-8b45-6daf77a83b61.md"
$mdFileUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc8cee7a51429c40490cdf627ee422436aad69ea/e2e/4f1cbecd-55f8-4daa-8b45-6daf77a83b61.md"

$zhXlf = "4f1cbecd-55f8-4daa-8b45-6daf77a83b61.84c8dc3deb0928ce5488d97451314089b4fea04c.zh-cn.xlf"
$deXlf = "4f1cbecd-55f8-4daa-8b45-6daf77a83b61.84c8dc3deb0928ce5488d97451314089b4fea04c.de-de.xlf"

$zhHandbackDate = "2016-08-21 19:05:12"
$deHandbackDate = "2016-08-21 19:05:18"

$status = "Handed back: in sync with en-US"

# Column width helper: ColumnWidth (characters) is quantized by Excel to
# 1/6-character steps, so pick the smallest input value that rounds up to
# the desired stored width.
$wideColWidth   = 39.0843333333333   # -> stored width 40
$mediumColWidth = 29.0843333333333   # -> stored width ~30 (closest to 29.9777047293527)

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $status
$wsZh.Range("C3").Value = $status

$wsZh.Range("I2").Value = $mdFile
$wsZh.Range("I3").Value = $mdFile
$wsZh.Range("I2:I3").Style = "HyperLink"

$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("J3").Value = $zhXlf

$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdFileUrl, "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdFileUrl, "", "", $mdFile)

$wsZh.Range("C3:C3").EntireColumn.ColumnWidth = $mediumColWidth
$wsZh.Range("I2:J3").EntireColumn.ColumnWidth = $wideColWidth

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $status
$wsDe.Range("C3").Value = $status

$wsDe.Range("I2").Value = $mdFile
$wsDe.Range("I3").Value = $mdFile
$wsDe.Range("I2:I3").Style = "HyperLink"

$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("J3").Value = $deXlf

$wsDe.Range("K2").Value = $deHandbackDate
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdFileUrl, "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdFileUrl, "", "", $mdFile)

$wsDe.Range("C3:C3").EntireColumn.ColumnWidth = $mediumColWidth
$wsDe.Range("I2:J3").EntireColumn.ColumnWidth = $wideColWidth

# ----------------------------------------------------------------------
# Overview sheet - zh-cn / de-de columns widen to match the new Status
# column width on the locale-specific sheets.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $mediumColWidth

Write-Host "Handback report generated"
